# "Add files via upload" - fills in the (previously empty) "Time estimation"
# column (D) on the "Sprint 2" sheet for each backlog item, and leaves the
# selection on E11, matching where the author's cursor ended up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "Sprint 2" is the sheet that was active/tabSelected

# Time estimations for rows 3-8 (column D), reusing "2h" where it repeats.
$ws.Range("D3").Value = "2h"
$ws.Range("D4").Value = "2.5h"
$ws.Range("D5").Value = "2h"
$ws.Range("D6").Value = "2h"
$ws.Range("D7").Value = "2.6h"
$ws.Range("D8").Value = "5h"

# Move/restore the selection to E11, as recorded in the saved view state.
$ws.Range("E11").Select() | Out-Null
